$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sprint1: User Story 1 (row 2) moves to "done" with a lower Est Time,
# and User Story 10 (row 11) gets a Status of "done" with an Est Time.
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("D2").Value = "done"
$sprint1.Range("F2").Value = 90
$sprint1.Range("D11").Value = "done"
$sprint1.Range("F11").Value = 120
$sprint1.Range("F11").Select()

# ---------------------------------------------------------------------
# Sprint3: just move the cursor (no data changes for this sheet).
# ---------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("Sprint3")
$sprint3.Range("A2").Select()

# ---------------------------------------------------------------------
# Stories: just move the cursor (no data changes for this sheet).
# ---------------------------------------------------------------------
$stories = $wb.Worksheets.Item("Stories")
$stories.Range("B17").Select()

# ---------------------------------------------------------------------
# Sprint2: two new backlog stories added, Sprint2 becomes the active tab.
# ---------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Activate()

$sprint2.Range("A2").Value = 15
$sprint2.Range("B2").Value = "Fewer than 15 siblings"
$sprint2.Range("C2").Value = "PD"
$sprint2.Range("D2").Value = "coding"

$sprint2.Range("A3").Value = 16
$sprint2.Range("B3").Value = "Male last names"
$sprint2.Range("C3").Value = "PD"
$sprint2.Range("D3").Value = "coding"

$sprint2.Columns.Item(2).ColumnWidth = 16.4

$sprint2.Range("B3").Select()
